$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.357.20'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.717.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4725'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2630'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06209'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.715.79'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07059'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('E12').Value = '  +3.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5904'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.418'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.07'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.346.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006815'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('E20').Value = '  +1.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.935.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.538'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.759'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.402'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '107.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.753'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.996'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.689'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07743'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04435'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.614'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9775'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6187'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9392'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '112.56'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +14.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.420'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.919'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01476'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.328'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +13.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3811'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1171'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.284'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05282'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.30'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.703'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.93%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3364'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.10%  '
